$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDA")

# Rename the shared string used for the little note in I9
$ws.Range("I9").Value = "Lookup alternative"

# The new formulas in column I/J push out the leftover blank placeholder
# cells that used to sit in column H for these two rows
$ws.Range("H11").ClearContents()
$ws.Range("H12").ClearContents()

# New spilled array formula =B3:B7 anchored at I11, spilling down to I15
$ws.Range("I11:I15").FormulaArray = "=B3:B7"

# Only the anchor cell I11 carries the same "date-ish" number format as the
# other leftover placeholder cells further down column A
$ws.Range("A15").Copy()
$ws.Range("I11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# One LOOKUP "alternative" formula per row, mirroring the existing
# INDEX/COUNTA block in J3:J7 but implemented with LOOKUP instead
$ws.Range("J11").FormulaArray = "=LOOKUP(2,1/(C3:G3<>""""),C3:G3)"
$ws.Range("J12").FormulaArray = "=LOOKUP(2,1/(C4:G4<>""""),C4:G4)"
$ws.Range("J13").FormulaArray = "=LOOKUP(2,1/(C5:G5<>""""),C5:G5)"
$ws.Range("J14").FormulaArray = "=LOOKUP(2,1/(C6:G6<>""""),C6:G6)"
$ws.Range("J15").FormulaArray = "=LOOKUP(2,1/(C7:G7<>""""),C7:G7)"

# Move the active selection to where the author left off
$ws.Range("F20").Select()
